# Pokemon Card Bright Tide list - update
# - add notes/questions to several card-count cells in column A (rows 23-29)
# - set column A width
# - update the scroll position / selection of the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update card-count text cells (append notes / markers) ---
$ws.Range("A23").Value = "2x 120 v"
$ws.Range("A24").Value = "1x 122 ? "
$ws.Range("A25").Value = "2x 123 ?"
$ws.Range("A26").Value = "2x 127 ? Heal ?"
$ws.Range("A27").Value = "1x 128 py script"
$ws.Range("A28").Value = "2x 134 ? "
$ws.Range("A29").Value = "1x 135 v"

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 15.5

# --- Sheet view: scroll position and active selection ---
$ws.Range("A26").Select()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
